# Applies the data updates described in the commit:
# "Elimna EC anteriores y se agregan nuevos, se modifica base de datos"
#
# Net effect on the "Hoja1" worksheet table (rows 16-20):
#   - Row 16 now holds LUCIANO AGUIRRE AGUDELO's data (was on row 20)
#   - Row 17 now holds TATIANA ZURIQUE DE ARCO's data (was on row 16)
#   - Row 18 stays NINI JOHANNA MOLINA GONZALEZ (unchanged)
#   - Row 19 stays YERINE GOMEZ SAENZ (unchanged)
#   - Row 20 now holds NELSON RAMIRO MONTOYA TORRES's data (was on row 17)
# Along with each worker's own "Valor Mora" (F) / "Salario Basico" (G) values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 16: Luciano Aguirre Agudelo
$ws.Range("C16").Value = "1053798208"
$ws.Range("D16").Value = "LUCIANO AGUIRRE AGUDELO"
$ws.Range("F16").Value = 80000
$ws.Range("G16").Value = 2000000

# Row 17: Tatiana Zurique De Arco
$ws.Range("C17").Value = "45549437"
$ws.Range("D17").Value = "TATIANA ZURIQUE DE ARCO"
$ws.Range("F17").Value = 24640
$ws.Range("G17").Value = 616000

# Row 20: Nelson Ramiro Montoya Torres
$ws.Range("C20").Value = "1090408750"
$ws.Range("D20").Value = "NELSON RAMIRO MONTOYA TORRES"
$ws.Range("F20").Value = 72000
$ws.Range("G20").Value = 1800000

$wb.Save()
